$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header labels (value/text relabeling, columns keep same position)
$ws.Range("A1").Value = "Número hogares"
$ws.Range("B1").Value = "Municipio 2ª residencia, nombre"
$ws.Range("C1").Value = "Aragón"
$ws.Range("D1").Value = "Municipio 2ª residencia, código"

# Row 3 - previously A3 was "null", now it mirrors "medida" (same as B3); C3 changes from "medida" to "dim"; D3 changes from "dim" to "null"
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"

# Row 4 - A4 "null" -> "xsd:int"; B4 "xsd:int" -> "xsd:string"; C4 "xsd:int" -> "URI-Comunidad"; D4 "URI-Comunidad" -> "null"
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "URI-Comunidad"
$ws.Range("D4").Value = "null"

$wb.Save()
